$d = $word.ActiveDocument
$sec = $d.Sections(1)

function Rename-InlineShapeInRange($rng, $newName) {
    $ishapes = $rng.InlineShapes
    for ($i = 1; $i -le $ishapes.Count; $i++) {
        $s = $ishapes.Item($i)
        $shp = $s.ConvertToShape()
        $shp.Name = $newName
        [void]$shp.ConvertToInlineShape()
    }
}

# Footers: Footers(1) is the "default" footer (footer2.xml, Pearson logo
# image2.png -> image1.png); Footers(2) is the "first page" footer
# (footer1.xml, also Pearson logo image2.png -> image1.png).
Rename-InlineShapeInRange $sec.Footers(1).Range "image1.png"
Rename-InlineShapeInRange $sec.Footers(2).Range "image1.png"

# Headers: Headers(1) is the "default" header (header2.xml, BTec logo
# image1.jpg -> image2.jpg); Headers(2) is the "first page" header
# (header1.xml, also BTec logo image1.jpg -> image2.jpg).
Rename-InlineShapeInRange $sec.Headers(1).Range "image2.jpg"
Rename-InlineShapeInRange $sec.Headers(2).Range "image2.jpg"
